# Auto-generated edit script: scheduled market-price refresh for Leve profit tables
# Updates currentAveragePrice/NQ/HQ and derived Leve profit columns (H-N) across all class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 21
$ws.Range("I20").Value = 21
$ws.Range("K20").Value = 21
$ws.Range("M20").Value = 209
$ws.Range("H32").Value = 35716236
$ws.Range("I32").Value = 2184
$ws.Range("K32").Value = 2184
$ws.Range("M32").Value = -1858
$ws.Range("H35").Value = 21
$ws.Range("I35").Value = 21
$ws.Range("K35").Value = 21
$ws.Range("M35").Value = 358
$ws.Range("H55").Value = 2449.6
$ws.Range("I55").Value = 562.25
$ws.Range("J55").Value = 9999
$ws.Range("K55").Value = 562.25
$ws.Range("L55").Value = 9999
$ws.Range("M55").Value = -348.25
$ws.Range("N55").Value = -10427
$ws.Range("H106").Value = 1899.1351
$ws.Range("I106").Value = 1610.1
$ws.Range("K106").Value = 1610.1
$ws.Range("M106").Value = -979.0999999999999
$ws.Range("H138").Value = 4521.098
$ws.Range("I138").Value = 6382.5
$ws.Range("J138").Value = 3613.0977
$ws.Range("K138").Value = 19147.5
$ws.Range("L138").Value = 10839.2931
$ws.Range("M138").Value = -14007.5
$ws.Range("N138").Value = -21119.2931

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 7666.778
$ws.Range("I13").Value = 9999
$ws.Range("J13").Value = 3002.3333
$ws.Range("K13").Value = 9999
$ws.Range("L13").Value = 3002.3333
$ws.Range("M13").Value = -9855
$ws.Range("N13").Value = -3290.3333
$ws.Range("H32").Value = 197638.22
$ws.Range("I32").Value = 223164.9
$ws.Range("J32").Value = 6188
$ws.Range("K32").Value = 223164.9
$ws.Range("L32").Value = 6188
$ws.Range("M32").Value = -222877.9
$ws.Range("N32").Value = -6762
$ws.Range("H61").Value = 921823.5600000001
$ws.Range("I61").Value = 2414.96
$ws.Range("K61").Value = 2414.96
$ws.Range("M61").Value = -2202.96
$ws.Range("H74").Value = 476112.06
$ws.Range("I74").Value = 1401.1613
$ws.Range("J74").Value = 1395864.5
$ws.Range("K74").Value = 1401.1613
$ws.Range("L74").Value = 1395864.5
$ws.Range("M74").Value = -527.1613
$ws.Range("N74").Value = -1397612.5
$ws.Range("H77").Value = 476112.06
$ws.Range("I77").Value = 1401.1613
$ws.Range("J77").Value = 1395864.5
$ws.Range("K77").Value = 7005.8065
$ws.Range("L77").Value = 6979322.5
$ws.Range("M77").Value = -2637.8065
$ws.Range("N77").Value = -6988058.5
$ws.Range("H132").Value = 2655.2778
$ws.Range("I132").Value = 2582.353
$ws.Range("K132").Value = 7747.059
$ws.Range("M132").Value = -5217.059
$ws.Range("H136").Value = 921823.5600000001
$ws.Range("I136").Value = 2414.96
$ws.Range("K136").Value = 7244.88
$ws.Range("M136").Value = -4694.88
$ws.Range("H139").Value = 98000
$ws.Range("J139").Value = 98000
$ws.Range("L139").Value = 98000
$ws.Range("N139").Value = -108280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 76923310
$ws.Range("J80").Value = 347.625
$ws.Range("L80").Value = 347.625
$ws.Range("N80").Value = -2343.625
$ws.Range("H83").Value = 76923310
$ws.Range("J83").Value = 347.625
$ws.Range("L83").Value = 1738.125
$ws.Range("N83").Value = -11722.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1560.56
$ws.Range("I134").Value = 1396.4286
$ws.Range("K134").Value = 4189.2858
$ws.Range("M134").Value = -1654.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 19995
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 19995
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 59985
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -60615
$ws.Range("H73").Value = 19995
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 19995
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 59985
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -62169

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 7328
$ws.Range("I10").Value = 3989
$ws.Range("J10").Value = 8997.5
$ws.Range("K10").Value = 3989
$ws.Range("L10").Value = 8997.5
$ws.Range("M10").Value = -3820
$ws.Range("N10").Value = -9335.5
$ws.Range("H29").Value = 999
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H35").Value = 13394.2
$ws.Range("J35").Value = 19000
$ws.Range("L35").Value = 19000
$ws.Range("N35").Value = -19596
$ws.Range("H41").Value = 7487.778
$ws.Range("I41").Value = 7487.778
$ws.Range("K41").Value = 7487.778
$ws.Range("M41").Value = -7132.778
$ws.Range("H80").Value = 12659789
$ws.Range("I80").Value = 165376.36
$ws.Range("K80").Value = 165376.36
$ws.Range("M80").Value = -164378.36
$ws.Range("H83").Value = 12659789
$ws.Range("I83").Value = 165376.36
$ws.Range("K83").Value = 826881.7999999999
$ws.Range("M83").Value = -821889.7999999999
$ws.Range("H126").Value = 2495
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1612590.8
$ws.Range("I132").Value = 17954
$ws.Range("K132").Value = 53862
$ws.Range("M132").Value = -51332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1058.3846
$ws.Range("I55").Value = 866.2308
$ws.Range("K55").Value = 866.2308
$ws.Range("M55").Value = -693.2308
$ws.Range("H132").Value = 3289.2593
$ws.Range("I132").Value = 2915.0625
$ws.Range("K132").Value = 8745.1875
$ws.Range("M132").Value = -6215.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 257974.75
$ws.Range("J3").Value = 343299.66
$ws.Range("L3").Value = 343299.66
$ws.Range("N3").Value = -343527.66
$ws.Range("H22").Value = 500
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -1086
$ws.Range("H23").Value = 2013.6666
$ws.Range("J23").Value = 2995.5
$ws.Range("L23").Value = 2995.5
$ws.Range("N23").Value = -3453.5
$ws.Range("H30").Value = 5000
$ws.Range("J30").Value = 5000
$ws.Range("L30").Value = 5000
$ws.Range("N30").Value = -5214
$ws.Range("H31").Value = 11320
$ws.Range("I31").Value = 12200
$ws.Range("K31").Value = 12200
$ws.Range("M31").Value = -11852
$ws.Range("H32").Value = 7845
$ws.Range("I32").Value = 7806.25
$ws.Range("K32").Value = 7806.25
$ws.Range("M32").Value = -7489.25
$ws.Range("H34").Value = 10000
$ws.Range("I34").Value = 10000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 10000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -9797
$ws.Range("N34").ClearContents()
$ws.Range("H113").Value = 734.4706
$ws.Range("J113").Value = 632
$ws.Range("L113").Value = 1896
$ws.Range("N113").Value = -6236
$ws.Range("H132").Value = 2859.6858
$ws.Range("I132").Value = 2565.7
$ws.Range("K132").Value = 7697.099999999999
$ws.Range("M132").Value = -5167.099999999999
$ws.Range("H136").Value = 24333.273
$ws.Range("I136").Value = 31842.25
$ws.Range("K136").Value = 95526.75
$ws.Range("M136").Value = -92976.75
